$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1275786666666667
$ws.Range("H2").Value = 0.382736
$ws.Range("I2").Value = 0.2993455218931061
$ws.Range("J2").Value = 0.2993455218931061
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04996866666666667
$ws.Range("N2").Value = 0.149906
$ws.Range("O2").Value = 0.06831051926220302
$ws.Range("P2").Value = 0.06831051926220301
$ws.Range("Q2").Value = 0.006374935868444445
$ws.Range("R2").Value = 0.057374422816
$ws.Range("S2").Value = 0.02044844803933324
$ws.Range("T2").Value = 0.02044844803933324

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1275786666666667
$ws.Range("H3").Value = 0.382736
$ws.Range("I3").Value = 0.2993455218931061
$ws.Range("J3").Value = 0.2993455218931061
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5898666666666667
$ws.Range("N3").Value = 1.7696
$ws.Range("O3").Value = 0.8063873019518528
$ws.Range("P3").Value = 0.8063873019518527
$ws.Range("Q3").Value = 0.07525440284444444
$ws.Range("R3").Value = 0.6772896256000001
$ws.Range("S3").Value = 0.2413884277507511
$ws.Range("T3").Value = 0.2413884277507511

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1275786666666667
$ws.Range("H4").Value = 0.382736
$ws.Range("I4").Value = 0.2993455218931061
$ws.Range("J4").Value = 0.2993455218931061
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.09165766666666668
$ws.Range("N4").Value = 0.274973
$ws.Range("O4").Value = 0.1253021787859442
$ws.Range("P4").Value = 0.1253021787859442
$ws.Range("Q4").Value = 0.01169356290311111
$ws.Range("R4").Value = 0.105242066128
$ws.Range("S4").Value = 0.03750864610302176
$ws.Range("T4").Value = 0.03750864610302175

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2986133333333333
$ws.Range("H5").Value = 0.89584
$ws.Range("I5").Value = 0.700654478106894
$ws.Range("J5").Value = 0.700654478106894
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04996866666666667
$ws.Range("N5").Value = 0.149906
$ws.Range("O5").Value = 0.06831051926220302
$ws.Range("P5").Value = 0.06831051926220301
$ws.Range("Q5").Value = 0.01492131011555556
$ws.Range("R5").Value = 0.13429179104
$ws.Range("S5").Value = 0.04786207122286978
$ws.Range("T5").Value = 0.04786207122286978

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2986133333333333
$ws.Range("H6").Value = 0.89584
$ws.Range("I6").Value = 0.700654478106894
$ws.Range("J6").Value = 0.700654478106894
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5898666666666667
$ws.Range("N6").Value = 1.7696
$ws.Range("O6").Value = 0.8063873019518528
$ws.Range("P6").Value = 0.8063873019518527
$ws.Range("Q6").Value = 0.1761420515555556
$ws.Range("R6").Value = 1.585278464
$ws.Range("S6").Value = 0.5649988742011017
$ws.Range("T6").Value = 0.5649988742011017

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2986133333333333
$ws.Range("H7").Value = 0.89584
$ws.Range("I7").Value = 0.700654478106894
$ws.Range("J7").Value = 0.700654478106894
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.09165766666666668
$ws.Range("N7").Value = 0.274973
$ws.Range("O7").Value = 0.1253021787859442
$ws.Range("P7").Value = 0.1253021787859442
$ws.Range("Q7").Value = 0.02737020136888889
$ws.Range("R7").Value = 0.24633181232
$ws.Range("S7").Value = 0.08779353268292246
$ws.Range("T7").Value = 0.08779353268292245
